{"js": "// Change the \"Heading 1\" paragraph's font size to 18pt (commit: \"Change H1 font to 18\").\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/style\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.style === \"Heading 1\") {\n    para.font.size = 18;\n    para.font.sizeBidirectional = 18;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Change the \"Heading 1\" paragraph's font size to 18pt (commit: \"Change H1 font to 18\").\n$d = $word.ActiveDocument\n\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Style.NameLocal -eq \"Heading 1\") {\n        $p.Range.Font.Size = 18\n        $p.Range.Font.SizeBi = 18\n    }\n}\n"}
